# Update "想去人数" (F column) figures for the latest scrape (gh-pages output
# generated at 456a3b4). The same underlying events appear on multiple
# sheets (展览, 演出, 全部类型) so the counts must be updated everywhere they
# occur.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# 展览 sheet
$wsExpo.Range("F2").Value = 612
$wsExpo.Range("F3").Value = 328
$wsExpo.Range("F8").Value = 1031
$wsExpo.Range("F9").Value = 3855
$wsExpo.Range("F10").Value = 76

# 演出 sheet
$wsShow.Range("F2").Value = 50

# 全部类型 sheet
$wsAll.Range("F2").Value = 612
$wsAll.Range("F3").Value = 328
$wsAll.Range("F8").Value = 1031
$wsAll.Range("F9").Value = 3855
$wsAll.Range("F10").Value = 76
$wsAll.Range("F11").Value = 50
